# CIERRE 24 DIC 2021
# Advance the payroll workbook from "SEMANA 51" (Dec 13-19, 2021) to
# "SEMANA 52" (Dec 20-26, 2021): update the week-header text, the per-employee
# days/amounts for this new week, and move the active selection back to B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Week header text (shared string), all formula-linked headers follow it ---
$ws.Range("B9").Value = "SEMANA   52  DEL    20      Al   26   DE   DICIEMBRE          2021"

# --- Block 1 (rows 3-7) updated figures ---
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 1467
$ws.Range("K4").Value = 867

# --- Block 2 (rows 21-26) updated figures ---
$ws.Range("K21").Value = 560

# --- Block 4 (rows 38-41) updated figures ---
$ws.Range("J39").Value = 3
$ws.Range("K39").Value = 1250

# --- Restore the active window to its resting scroll/selection state ---
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
